# Append new listing rows (341-345) to the carousell listings sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("1313409423", "https://www.carousell.sg/p/1313409423", "little_stiches", "2 days ago", "Apple Pencil Gen 1 (Free case)", "S$80", "Like new", "https://media.karousell.com/media/photos/products/2024/7/9/apple_pencil_gen_1_free_case_1720547160_018d4006_progressive_thumbnail.jpg"),
    @("1313128588", "https://www.carousell.sg/p/1313128588", "zaw_m_h", "3 days ago", "Apple Pencil Gen 1", "S$110", "Lightly used", "https://media.karousell.com/media/photos/products/2024/7/8/apple_pencil_gen_1_1720430285_84c0e654_progressive_thumbnail.jpg"),
    @("1312992098", "https://www.carousell.sg/p/1312992098", "heartaimer", "4 days ago", "Apple Pencil Gen 1", "S$75", "Like new", "https://media.karousell.com/media/photos/products/2024/7/7/apple_pencil_gen_1_1720362436_2b2d690f_progressive_thumbnail.jpg"),
    @("1312916504", "https://www.carousell.sg/p/1312916504", "auzra35", "4 days ago", "Apple Pencil Gen 1", "S$90", "Like new", "https://media.karousell.com/media/photos/products/2024/7/7/apple_pencil_gen_1_1720340083_9423c6b1_progressive_thumbnail.jpg"),
    @("1312513290", "https://www.carousell.sg/p/1312513290", "yljc199702018", "6 days ago", "Apple Pencil Gen 1", "S$100", "Like new", "https://media.karousell.com/media/photos/products/2024/7/5/apple_pencil_gen_1_1720157937_2fe922b0_progressive_thumbnail.jpg")
)

$startRow = 341
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le 8; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c - 1]
    }
}
